# Update "想去人数" (number of attendees) figures that changed between
# the two data pulls for the 南宁漫展信息 workbook.
#
# Sheet "展览" (exhibitions):
#   F3: 1130 -> 1134
#   F4: 2569 -> 2578
#   F5: 222  -> 223
#
# Sheet "全部类型" (all types, aggregated view of the same events):
#   F5: 1130 -> 1134
#   F6: 2569 -> 2578
#   F8: 222  -> 223

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 1134
$wsExhibition.Range("F4").Value = 2578
$wsExhibition.Range("F5").Value = 223

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 1134
$wsAll.Range("F6").Value = 2578
$wsAll.Range("F8").Value = 223
